$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("J3").Value = "dim"
$ws.Range("J4").Value = "URI-Municipio"
